$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial number for every data row
# (rows 2 through 312). Update the value from 45206 (2023-10-07) to
# 45208 (2023-10-09) for all of them.
$ws.Range("C2:C312").Value = 45208
